# Update the LR-pairs data sheet with refreshed TPM values.
# New sending cluster "ECs" is added (rows 2-5), and the previously existing
# FAPs / Resolving-Mac blocks shift down to rows 6-9 / 10-13 with updated values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs", "Lta", "Tnfrsf14", "ECs", 1, 0.3333333333333333, 0.103879, 0.311637, 0.1899090787212519, 0.1899090787212519, 3, 1, 5.209944, 15.629832, 0.1751928672265232, 0.1751928672265232, 0.541203772776, 4.870833954984001, 0.03327071601352363, 0.03327071601352363),
    @("ECs", "Lta", "Tnfrsf14", "FAPs", 1, 0.3333333333333333, 0.103879, 0.311637, 0.1899090787212519, 0.1899090787212519, 3, 1, 6.497702, 19.493106, 0.2184958310038485, 0.2184958310038485, 0.674974786058, 6.074773074522001, 0.04149434197037521, 0.04149434197037521),
    @("ECs", "Lta", "Tnfrsf14", "MuSCs", 1, 0.3333333333333333, 0.103879, 0.311637, 0.1899090787212519, 0.1899090787212519, 3, 1, 1.714656666666667, 5.14397, 0.057658127945791, 0.057658127945791, 0.1781168198766667, 1.60305137889, 0.01094980195897724, 0.01094980195897724),
    @("ECs", "Lta", "Tnfrsf14", "Resolving-Mac", 1, 0.3333333333333333, 0.103879, 0.311637, 0.1899090787212519, 0.1899090787212519, 3, 1, 16.31603133333333, 48.948094, 0.5486531738238374, 0.5486531738238374, 1.694893018875333, 15.254037169878, 0.1041942187783759, 0.1041942187783759),
    @("FAPs", "Lta", "Tnfrsf14", "ECs", 2, 0.6666666666666666, 0.257314, 0.771942, 0.4704152396738534, 0.4704152396738534, 3, 1, 5.209944, 15.629832, 0.1751928672265232, 0.1751928672265232, 1.340591530416, 12.065323773744, 0.08241339462551449, 0.08241339462551449),
    @("FAPs", "Lta", "Tnfrsf14", "FAPs", 2, 0.6666666666666666, 0.257314, 0.771942, 0.4704152396738534, 0.4704152396738534, 3, 1, 6.497702, 19.493106, 0.2184958310038485, 0.2184958310038485, 1.671949692428, 15.047547231852, 0.1027837687094131, 0.1027837687094131),
    @("FAPs", "Lta", "Tnfrsf14", "MuSCs", 2, 0.6666666666666666, 0.257314, 0.771942, 0.4704152396738534, 0.4704152396738534, 3, 1, 1.714656666666667, 5.14397, 0.057658127945791, 0.057658127945791, 0.4412051655266667, 3.97084648974, 0.02712326207676498, 0.02712326207676498),
    @("FAPs", "Lta", "Tnfrsf14", "Resolving-Mac", 2, 0.6666666666666666, 0.257314, 0.771942, 0.4704152396738534, 0.4704152396738534, 3, 1, 16.31603133333333, 48.948094, 0.5486531738238374, 0.5486531738238374, 4.198343286505333, 37.785089578548, 0.2580948142621608, 0.2580948142621608),
    @("Resolving-Mac", "Lta", "Tnfrsf14", "ECs", 2, 0.6666666666666666, 0.1858003333333333, 0.557401, 0.3396756816048946, 0.3396756816048946, 3, 1, 5.209944, 15.629832, 0.1751928672265232, 0.1751928672265232, 0.9680093318480001, 8.712083986632, 0.05950875658748507, 0.05950875658748507),
    @("Resolving-Mac", "Lta", "Tnfrsf14", "FAPs", 2, 0.6666666666666666, 0.1858003333333333, 0.557401, 0.3396756816048946, 0.3396756816048946, 3, 1, 6.497702, 19.493106, 0.2184958310038485, 0.2184958310038485, 1.207275197500667, 10.865476777506, 0.0742177203240601, 0.0742177203240601),
    @("Resolving-Mac", "Lta", "Tnfrsf14", "MuSCs", 2, 0.6666666666666666, 0.1858003333333333, 0.557401, 0.3396756816048946, 0.3396756816048946, 3, 1, 1.714656666666667, 5.14397, 0.057658127945791, 0.057658127945791, 0.3185837802188889, 2.86725402197, 0.01958506391004878, 0.01958506391004878),
    @("Resolving-Mac", "Lta", "Tnfrsf14", "Resolving-Mac", 2, 0.6666666666666666, 0.1858003333333333, 0.557401, 0.3396756816048946, 0.3396756816048946, 3, 1, 16.31603133333333, 48.948094, 0.5486531738238374, 0.5486531738238374, 3.031524060410444, 27.283716543694, 0.1863641407833007, 0.1863641407833007)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    for ($col = 1; $col -le $values.Count; $col++) {
        $ws.Cells.Item($row, $col).Value = $values[$col - 1]
    }
}
